# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 267
$ws1.Range("F3").Value = 1379
$ws1.Range("F4").Value = 159
$ws1.Range("F6").Value = 234
$ws1.Range("F10").Value = 136
$ws1.Range("F11").Value = 4633
$ws1.Range("F12").Value = 6896
$ws1.Range("F18").Value = 4142
$ws1.Range("F19").Value = 733
$ws1.Range("F22").Value = 2730
$ws1.Range("F25").Value = 170
$ws1.Range("F26").Value = 369
$ws1.Range("F28").Value = 401
$ws1.Range("F31").Value = 1634
$ws1.Range("F32").Value = 1028
$ws1.Range("F34").Value = 364
$ws1.Range("F36").Value = 549
$ws1.Range("F40").Value = 165

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 267
$ws4.Range("F3").Value = 1379
$ws4.Range("F4").Value = 159
$ws4.Range("F6").Value = 234
$ws4.Range("F10").Value = 136
$ws4.Range("F11").Value = 4633
$ws4.Range("F12").Value = 6896
$ws4.Range("F18").Value = 4142
$ws4.Range("F19").Value = 734
$ws4.Range("F22").Value = 2730
$ws4.Range("F25").Value = 170
$ws4.Range("F26").Value = 369
$ws4.Range("F28").Value = 401
$ws4.Range("F31").Value = 1634
$ws4.Range("F32").Value = 1028
$ws4.Range("F34").Value = 364
$ws4.Range("F36").Value = 549
$ws4.Range("F40").Value = 165
